$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet: "VENTAS POR GRUPO" (columns A:R) - insert new client row
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Rows.Item(18).Insert()
$ws1.Cells.Item(18,1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(18,2).Value = "CORONADO MONTERO LIDA VERONICA"
for ($c = 3; $c -le 18; $c++) {
  $ws1.Cells.Item(18,$c).Value = 0
}
# the summary row (was row 52, now shifted to row 53) mentions "X de 50"; bump it to "X de 51"
for ($c = 3; $c -le 18; $c++) {
  $cell = $ws1.Cells.Item(53, $c)
  $cell.Value = ($cell.Value() -replace "de 50", "de 51")
}

# ---------------------------------------------------------------
# Sheet: "VENTA MENSUAL" (columns A:G) - same new client row
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Rows.Item(18).Insert()
$ws2.Cells.Item(18,1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(18,2).Value = "CORONADO MONTERO LIDA VERONICA"
for ($c = 3; $c -le 7; $c++) {
  $ws2.Cells.Item(18,$c).Value = 0
}

Write-Host "edit complete"
